$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.01942288123009
$ws.Range("C2").Value = 5.162034774222009
$ws.Range("E2").Value = 22.0718331648662
$ws.Range("F2").Value = 42.0137265407898
$ws.Range("G2").Value = 32.64878186171941
$ws.Range("H2").Value = 15.27365411732807
$ws.Range("I2").Value = 21.72367760119745
$ws.Range("J2").Value = 8.298844723406516
$ws.Range("K2").Value = 10.16049146873876
$ws.Range("N2").Value = 18.60313205583817
$ws.Range("B3").Value = 9.73684480672199
$ws.Range("C3").Value = 4.945389228703379
$ws.Range("E3").Value = 21.77519265401421
$ws.Range("F3").Value = 41.81531558138956
$ws.Range("G3").Value = 32.71401937206285
$ws.Range("H3").Value = 15.32688032833076
$ws.Range("I3").Value = 21.81785061047033
$ws.Range("J3").Value = 8.323063595870696
$ws.Range("K3").Value = 9.973074050693921
$ws.Range("N3").Value = 18.66818125429065
$ws.Range("B4").Value = 9.561124857280868
$ws.Range("C4").Value = 4.808526192994369
$ws.Range("E4").Value = 21.59637723113245
$ws.Range("F4").Value = 41.70492786115938
$ws.Range("G4").Value = 32.76520126959029
$ws.Range("H4").Value = 15.36226044584568
$ws.Range("I4").Value = 21.88008614636295
$ws.Range("J4").Value = 8.338880304591587
$ws.Range("K4").Value = 9.858270103827671
$ws.Range("N4").Value = 18.70994175724095
$ws.Range("B5").Value = 9.489075499814021
$ws.Range("C5").Value = 4.751875869867916
$ws.Range("E5").Value = 21.52442960561703
$ws.Range("F5").Value = 41.66285305243019
$ws.Range("G5").Value = 32.78884150136205
$ws.Range("H5").Value = 15.37735606323586
$ws.Range("I5").Value = 21.90655550965254
$ws.Range("J5").Value = 8.345564044420261
$ws.Range("K5").Value = 9.811618874906642
$ws.Range("N5").Value = 18.72741864533266
$ws.Range("B6").Value = 9.477088554206778
$ws.Range("C6").Value = 4.742418922943798
$ws.Range("E6").Value = 21.51254085360105
$ws.Range("F6").Value = 41.65604314136691
$ws.Range("G6").Value = 32.79293459138968
$ws.Range("H6").Value = 15.37990360245487
$ws.Range("I6").Value = 21.91101758867519
$ws.Range("J6").Value = 8.346688278109893
$ws.Range("K6").Value = 9.803882355197207
$ws.Range("N6").Value = 18.73034844403187
$ws.Range("B7").Value = 9.560154807300139
$ws.Range("C7").Value = 4.807765615193698
$ws.Range("E7").Value = 21.59540307912419
$ws.Range("F7").Value = 41.70434860659501
$ws.Range("G7").Value = 32.76550884123044
$ws.Range("H7").Value = 15.36246128671902
$ws.Range("I7").Value = 21.88043863807557
$ws.Range("J7").Value = 8.338969478475493
$ws.Range("K7").Value = 9.85764032799095
$ws.Range("N7").Value = 18.71017559565957
$ws.Range("B8").Value = 9.922516854584224
$ws.Range("C8").Value = 5.088182616122805
$ws.Range("E8").Value = 21.96891360847002
$ws.Range("F8").Value = 41.94296179074494
$ws.Range("G8").Value = 32.66895830690954
$ws.Range("H8").Value = 15.29144588192652
$ws.Range("I8").Value = 21.75523147027434
$ws.Range("J8").Value = 8.306999222439455
$ws.Range("K8").Value = 10.09585186660601
$ws.Range("N8").Value = 18.6251841682634
$ws.Range("B9").Value = 10.61082047540656
$ws.Range("C9").Value = 5.604229130114512
$ws.Range("E9").Value = 22.72391067074292
$ws.Range("F9").Value = 42.49997267441973
$ws.Range("G9").Value = 32.56848379058581
$ws.Range("H9").Value = 15.17363420053705
$ws.Range("I9").Value = 21.54479590445291
$ws.Range("J9").Value = 8.251796939994252
$ws.Range("K9").Value = 10.56237360579627
$ws.Range("N9").Value = 18.4728885889533
$ws.Range("B10").Value = 11.09715626410232
$ws.Range("C10").Value = 5.95889509072594
$ws.Range("E10").Value = 23.28719074778549
$ws.Range("F10").Value = 42.96101358049602
$ws.Range("G10").Value = 32.54953489687604
$ws.Range("H10").Value = 15.10020058167122
$ws.Range("I10").Value = 21.41169556254528
$ws.Range("J10").Value = 8.215783812808281
$ws.Range("K10").Value = 10.90117223949094
$ws.Range("N10").Value = 18.36966086494122
$ws.Range("B11").Value = 11.31313305988461
$ws.Range("C11").Value = 6.114289416212258
$ws.Range("E11").Value = 23.5442166537984
$ws.Range("F11").Value = 43.18138655123131
$ws.Range("G11").Value = 32.55294405686001
$ws.Range("H11").Value = 15.06965333088381
$ws.Range("I11").Value = 21.35583956449878
$ws.Range("J11").Value = 8.200382445507604
$ws.Range("K11").Value = 11.05372615268516
$ws.Range("N11").Value = 18.32456022631229
$ws.Range("B12").Value = 11.39408317380603
$ws.Range("C12").Value = 6.17223508481991
$ws.Range("E12").Value = 23.64156958037253
$ws.Range("F12").Value = 43.26630941024339
$ws.Range("G12").Value = 32.55597134295021
$ws.Range("H12").Value = 15.05849775729234
$ws.Range("I12").Value = 21.33536526747704
$ws.Range("J12").Value = 8.194691101699075
$ws.Range("K12").Value = 11.11121289700562
$ws.Range("N12").Value = 18.30774743508411
$ws.Range("B13").Value = 11.3766874643171
$ws.Range("C13").Value = 6.159796052536213
$ws.Range("E13").Value = 23.62060328818064
$ws.Range("F13").Value = 43.24795522456126
$ws.Range("G13").Value = 32.55524205540566
$ws.Range("H13").Value = 15.06088196926249
$ws.Range("I13").Value = 21.33974461747804
$ws.Range("J13").Value = 8.195910575692372
$ws.Range("K13").Value = 11.09884551653727
$ws.Range("N13").Value = 18.31135656866551
$ws.Range("B14").Value = 11.31981004854273
$ws.Range("C14").Value = 6.119074882733493
$ws.Range("E14").Value = 23.55222599092395
$ws.Range("F14").Value = 43.18834404032219
$ws.Range("G14").Value = 32.55315828437843
$ws.Range("H14").Value = 15.06872729204321
$ws.Range("I14").Value = 21.35414154634696
$ws.Range("J14").Value = 8.199911394527506
$ws.Range("K14").Value = 11.05846159002042
$ws.Range("N14").Value = 18.32317170775329
$ws.Range("B15").Value = 11.28485994447805
$ws.Range("C15").Value = 6.09401374974779
$ws.Range("E15").Value = 23.51034324422445
$ws.Range("F15").Value = 43.15202044283294
$ws.Range("G15").Value = 32.5521081959118
$ws.Range("H15").Value = 15.07358646458984
$ws.Range("I15").Value = 21.36304834023729
$ws.Range("J15").Value = 8.202380342494051
$ws.Range("K15").Value = 11.0336868967858
$ws.Range("N15").Value = 18.33044340183322
$ws.Range("B16").Value = 11.08292888347795
$ws.Range("C16").Value = 5.94861617313753
$ws.Range("E16").Value = 23.27040156470791
$ws.Range("F16").Value = 42.94682096235001
$ws.Range("G16").Value = 32.54955475786439
$ws.Range("H16").Value = 15.10225451521509
$ws.Range("I16").Value = 21.41544050736469
$ws.Range("J16").Value = 8.216810048769752
$ws.Range("K16").Value = 10.89116602690978
$ws.Range("N16").Value = 18.37264557232061
$ws.Range("B17").Value = 10.95764474773898
$ws.Range("C17").Value = 5.857863445620245
$ws.Range("E17").Value = 23.12334175403753
$ws.Range("F17").Value = 42.82362337073779
$ws.Range("G17").Value = 32.55107448700998
$ws.Range("H17").Value = 15.12057418159171
$ws.Range("I17").Value = 21.4487848743118
$ws.Range("J17").Value = 8.225913314043504
$ws.Range("K17").Value = 10.80329136697044
$ws.Range("N17").Value = 18.39901018756629
$ws.Range("B18").Value = 10.88509379739326
$ws.Range("C18").Value = 5.805107481374138
$ws.Range("E18").Value = 23.03883621035666
$ws.Range("F18").Value = 42.75376884250888
$ws.Range("G18").Value = 32.55308055566643
$ws.Range("H18").Value = 15.13138012306817
$ws.Range("I18").Value = 21.46840515774837
$ws.Range("J18").Value = 8.231241641056014
$ws.Range("K18").Value = 10.75260387845172
$ws.Range("N18").Value = 18.41434939739635
$ws.Range("B19").Value = 10.86044753146618
$ws.Range("C19").Value = 5.787150932044419
$ws.Range("E19").Value = 23.01024050348905
$ws.Range("F19").Value = 42.730291708096
$ws.Range("G19").Value = 32.55395396200149
$ws.Range("H19").Value = 15.13508499124729
$ws.Range("I19").Value = 21.4751240001499
$ws.Range("J19").Value = 8.233061594064628
$ws.Range("K19").Value = 10.73541895789092
$ws.Range("N19").Value = 18.41957308500315
$ws.Range("B20").Value = 10.971032855003
$ws.Range("C20").Value = 5.867582255668021
$ws.Range("E20").Value = 23.13898896859067
$ws.Range("F20").Value = 42.83663431331343
$ws.Range("G20").Value = 32.55079550225349
$ws.Range("H20").Value = 15.11859617870018
$ws.Range("I20").Value = 21.4451895997658
$ws.Range("J20").Value = 8.224934698012822
$ws.Range("K20").Value = 10.81266113031931
$ws.Range("N20").Value = 18.39618552931373
$ws.Range("B21").Value = 11.33653959678615
$ws.Range("C21").Value = 6.131060388005948
$ws.Range("E21").Value = 23.57231015168488
$ws.Range("F21").Value = 43.20581380751229
$ws.Range("G21").Value = 32.55372317251852
$ws.Range("H21").Value = 15.0664117406699
$ws.Range("I21").Value = 21.34989442177257
$ws.Range("J21").Value = 8.198732437386221
$ws.Range("K21").Value = 11.07033142191271
$ws.Range("N21").Value = 18.31969411092352
$ws.Range("B22").Value = 11.57051729498457
$ws.Range("C22").Value = 6.298003807779954
$ws.Range("E22").Value = 23.8556031356265
$ws.Range("F22").Value = 43.45564721144826
$ws.Range("G22").Value = 32.56575884085561
$ws.Range("H22").Value = 15.03470816183383
$ws.Range("I22").Value = 21.29156166881283
$ws.Range("J22").Value = 8.182428429414323
$ws.Range("K22").Value = 11.23706531638883
$ws.Range("N22").Value = 18.27125145243246
$ws.Range("B23").Value = 11.4461117845039
$ws.Range("C23").Value = 6.2093966344286
$ws.Range("E23").Value = 23.70442515603966
$ws.Range("F23").Value = 43.32154374870239
$ws.Range("G23").Value = 32.55840731054299
$ws.Range("H23").Value = 15.05140884093772
$ws.Range("I23").Value = 21.32233288079908
$ws.Range("J23").Value = 8.191055180913073
$ws.Range("K23").Value = 11.14824693379983
$ws.Range("N23").Value = 18.29696493304883
$ws.Range("B24").Value = 10.96498171595813
$ws.Range("C24").Value = 5.863190187966486
$ws.Range("E24").Value = 23.1319147279392
$ws.Range("F24").Value = 42.83074902957087
$ws.Range("G24").Value = 32.5509181044622
$ws.Range("H24").Value = 15.11948958121703
$ws.Range("I24").Value = 21.44681362128477
$ws.Range("J24").Value = 8.225376835194046
$ws.Range("K24").Value = 10.80842557664693
$ws.Range("N24").Value = 18.3974619907451
$ws.Range("B25").Value = 10.4276282874196
$ws.Range("C25").Value = 5.468679752806735
$ws.Range("E25").Value = 22.5177834662112
$ws.Range("F25").Value = 42.33999153681492
$ws.Range("G25").Value = 32.58607320245023
$ws.Range("H25").Value = 15.20320429745275
$ws.Range("I25").Value = 21.59795573499938
$ws.Range("J25").Value = 8.265931022825603
$ws.Range("K25").Value = 10.43661029468531
$ws.Range("N25").Value = 18.51256008078095
